$d = $word.ActiveDocument

# Locate the final paragraph of the document, which currently reads
# "[histograms] [outliers] [pair plots] [correlation]" and has style
# FirstParagraph. We replace its whole range (text + paragraph mark)
# with: a new narrative paragraph describing the `train` dataset
# statistics, a "Table 3" caption, the new statistics table, and then
# re-insert the original sentence as its own paragraph now styled as
# BodyText (matching the target edit).
$count = $d.Paragraphs.Count
$target = $d.Paragraphs($count)

if ($target.Range.Text.TrimEnd([char]13) -ne "[histograms] [outliers] [pair plots] [correlation]") {
    throw "Unexpected final paragraph text: " + $target.Range.Text
}

$rng = $target.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:pStyle w:val="FirstParagraph"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Throughout the EDA and model training, we will only use the</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="VerbatimChar"/>
        </w:rPr>
        <w:t xml:space="preserve">train</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">dataset which has a total of 1225 rows, with a minimum value</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">of -0.29562 (X19) and a</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">maximum value of 16.4051897</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">(X5). The variables have different scales</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">and variances:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="TableCaption"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Table 3: Statistics across Variables</w:t>
      </w:r>
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="Table"/>
        <w:tblW w:type="auto" w:w="0"/>
        <w:tblLook w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0" w:val="0020"/>
        <w:jc w:val="start"/>
        <w:tblCaption w:val="Table 3: Statistics across Variables"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="2640"/>
        <w:gridCol w:w="2640"/>
        <w:gridCol w:w="2640"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:tblHeader w:val="true"/>
        </w:trPr>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Statistic</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Min</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Max</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Mean</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">0.2425769 (X8)</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">13.8458721 (X5)</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Variance</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">0.033448 (X8)</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">1.5341953 (X4)</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Range</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">1.0023615 (X8)</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Compact"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">8.6214209 (X14)</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">[histograms] [outliers] [pair plots] [correlation]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)

Write-Host "Paragraphs after insert:" $d.Paragraphs.Count
Write-Host "Tables after insert:" $d.Tables.Count
